$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2612.5
$ws.Range("I18").Value = 1725
$ws.Range("J18").Value = 3500
$ws.Range("K18").Value = 1725
$ws.Range("L18").Value = 3500
$ws.Range("M18").Value = -1441
$ws.Range("N18").Value = -4068
$ws.Range("H58").Value = 1512.7858
$ws.Range("I58").Value = 937.9
$ws.Range("J58").Value = 2950
$ws.Range("K58").Value = 2813.7
$ws.Range("L58").Value = 8850
$ws.Range("M58").Value = -2663.7
$ws.Range("N58").Value = -9150
$ws.Range("H98").Value = 1945.4166
$ws.Range("I98").Value = 1934.5
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 1934.5
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -436.5
$ws.Range("N98").Value = -4996
$ws.Range("H122").Value = 1945.4166
$ws.Range("I122").Value = 1934.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5803.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3353.5
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 1822.0227
$ws.Range("I132").Value = 1903.5278
$ws.Range("J132").Value = 1455.25
$ws.Range("K132").Value = 5710.5834
$ws.Range("L132").Value = 4365.75
$ws.Range("M132").Value = -3180.5834
$ws.Range("N132").Value = -9425.75
$ws.Range("H137").Value = 1015
$ws.Range("I137").Value = 990
$ws.Range("K137").Value = 2970
$ws.Range("M137").Value = -420

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 274.33334
$ws.Range("I5").Value = 111.5
$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 111.5
$ws.Range("L5").Value = 600
$ws.Range("M5").Value = 0.5
$ws.Range("N5").Value = -824

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 274.33334
$ws.Range("I4").Value = 111.5
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 111.5
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = 3.5
$ws.Range("N4").Value = -830

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 3250000
$ws.Range("I6").Value = 5000000
$ws.Range("J6").Value = 1500000
$ws.Range("K6").Value = 5000000
$ws.Range("L6").Value = 1500000
$ws.Range("M6").Value = -4999887
$ws.Range("N6").Value = -1500226
$ws.Range("H7").Value = 321.27274
$ws.Range("I7").Value = 392.85715
$ws.Range("K7").Value = 392.85715
$ws.Range("M7").Value = -279.85715
$ws.Range("H17").Value = 500
$ws.Range("J17").Value = 500
$ws.Range("L17").Value = 500
$ws.Range("N17").Value = -848
$ws.Range("I22").Value = 1517.5555
$ws.Range("J22").Value = 666.6667
$ws.Range("K22").Value = 1517.5555
$ws.Range("L22").Value = 666.6667
$ws.Range("M22").Value = -1167.5555
$ws.Range("N22").Value = -1366.6667
$ws.Range("H41").Value = 11428
$ws.Range("I41").Value = 4000
$ws.Range("K41").Value = 4000
$ws.Range("M41").Value = -3572
$ws.Range("H50").Value = 9187.25
$ws.Range("J50").Value = 9187.25
$ws.Range("L50").Value = 9187.25
$ws.Range("N50").Value = -10437.25
$ws.Range("H51").Value = 8733.333000000001
$ws.Range("J51").Value = 9386.666999999999
$ws.Range("L51").Value = 9386.666999999999
$ws.Range("N51").Value = -10858.667
$ws.Range("H59").Value = 14881
$ws.Range("I59").Value = 8000
$ws.Range("J59").Value = 16601.25
$ws.Range("K59").Value = 8000
$ws.Range("L59").Value = 16601.25
$ws.Range("M59").Value = -6855
$ws.Range("N59").Value = -18891.25
$ws.Range("H60").Value = 6984.6
$ws.Range("I60").Value = 5050
$ws.Range("J60").Value = 8274.333000000001
$ws.Range("K60").Value = 5050
$ws.Range("L60").Value = 8274.333000000001
$ws.Range("M60").Value = -4539
$ws.Range("N60").Value = -9296.333000000001
$ws.Range("H61").Value = 8733.333000000001
$ws.Range("J61").Value = 9386.666999999999
$ws.Range("L61").Value = 9386.666999999999
$ws.Range("N61").Value = -10082.667
$ws.Range("H68").Value = 15661.875
$ws.Range("J68").Value = 15661.875
$ws.Range("L68").Value = 15661.875
$ws.Range("N68").Value = -17159.875
$ws.Range("H71").Value = 15661.875
$ws.Range("J71").Value = 15661.875
$ws.Range("L71").Value = 46985.625
$ws.Range("N71").Value = -54473.625
$ws.Range("H74").Value = 13679.5
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 13679.5
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 13679.5
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -15427.5
$ws.Range("H77").Value = 13679.5
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 13679.5
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 41038.5
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -49774.5
$ws.Range("H132").Value = 1731.2858
$ws.Range("I132").Value = 1116.3125
$ws.Range("K132").Value = 3348.9375
$ws.Range("M132").Value = -818.9375
$ws.Range("H134").Value = 62501336
$ws.Range("I134").Value = 1616.6666
$ws.Range("J134").Value = 250000500
$ws.Range("K134").Value = 4849.9998
$ws.Range("L134").Value = 750001500
$ws.Range("M134").Value = -2314.9998
$ws.Range("N134").Value = -750006570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3958.65
$ws.Range("I70").Value = 4372.364
$ws.Range("J70").Value = 3453
$ws.Range("K70").Value = 4372.364
$ws.Range("L70").Value = 3453
$ws.Range("M70").Value = -4102.364
$ws.Range("N70").Value = -3993
$ws.Range("H73").Value = 3958.65
$ws.Range("I73").Value = 4372.364
$ws.Range("J73").Value = 3453
$ws.Range("K73").Value = 4372.364
$ws.Range("L73").Value = 3453
$ws.Range("M73").Value = -3436.364
$ws.Range("N73").Value = -5325

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 14964
$ws.Range("I122").Value = 19018.666
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 57055.99800000001
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -54605.99800000001
$ws.Range("N122").Value = -13300
$ws.Range("H132").Value = 2944.5862
$ws.Range("I132").Value = 2580
$ws.Range("K132").Value = 7740
$ws.Range("M132").Value = -5210

Write-Host "Applied all cell updates"